# Section 5 & 6 reference renumbering
# [19] -> [24], [11] -> [10], [20] -> [10] (first occurrence),
# [20] -> [25] (second/third occurrences), [21] -> [26], [22] -> [27]
#
# Each Find/Replace is scoped to the specific paragraph that holds the
# citation so the many re-used numbers ("19", "20", "21", "22", "11", "9", …)
# elsewhere in the document are left untouched.

$d = $word.ActiveDocument

# Every Find is scoped to the one paragraph that holds the citation (each
# "[NN]" string is verified unique within that paragraph) and the matched
# text is exactly the "[NN]" bracket run itself — never any surrounding
# plain-formatted words — so the highlight/color run-formatting carried by
# the bracket is preserved on the replacement instead of being overwritten
# by the (differently formatted) neighbouring text's style.

# --- 5.1 Suicide Prevention Strategy -------------------------------------

# "...implementation of national strategies" [19]." -> "[24]"
$r = $d.Paragraphs.Item(7).Range
$r.Find.Execute("[19]", $false, $false, $false, $false, $false, $true, 1, $false, "[24]", 2) | Out-Null

# "...study in 2018 [11] which contained..." -> "[10]"
$r = $d.Paragraphs.Item(8).Range
$r.Find.Execute("[11]", $false, $false, $false, $false, $false, $true, 1, $false, "[10]", 2) | Out-Null

# "...adopted by the government" [20]." -> "[10]"
$r = $d.Paragraphs.Item(8).Range
$r.Find.Execute("[20]", $false, $false, $false, $false, $false, $true, 1, $false, "[10]", 2) | Out-Null

# "...success in reducing suicide [20]:" -> "[25]"
$r = $d.Paragraphs.Item(12).Range
$r.Find.Execute("[20]", $false, $false, $false, $false, $false, $true, 1, $false, "[25]", 2) | Out-Null

# Reference list entry: "[19] https://www.suicideinfo.ca/..." -> "[24]"
$r = $d.Paragraphs.Item(23).Range
$r.Find.Execute("[19]", $false, $false, $false, $false, $false, $true, 1, $false, "[24]", 2) | Out-Null

# Reference list entry: "[20] https://www.who.int/.../world_report_2014/en/" -> "[25]"
$r = $d.Paragraphs.Item(26).Range
$r.Find.Execute("[20]", $false, $false, $false, $false, $false, $true, 1, $false, "[25]", 2) | Out-Null

# --- 5.2 Alcohol Intake ----------------------------------------------------

# "...decrease in suicide. [21]" -> "[26]"
$r = $d.Paragraphs.Item(30).Range
$r.Find.Execute("[21]", $false, $false, $false, $false, $false, $true, 1, $false, "[26]", 2) | Out-Null

# Reference list entry: "[21] https://academic.oup.com/..." -> "[26]"
$r = $d.Paragraphs.Item(33).Range
$r.Find.Execute("[21]", $false, $false, $false, $false, $false, $true, 1, $false, "[26]", 2) | Out-Null

# --- 5.3 GDP Per Capita ------------------------------------------------

# "...other associated industries. [22] With these..." -> "[27]"
$r = $d.Paragraphs.Item(36).Range
$r.Find.Execute("[22]", $false, $false, $false, $false, $false, $true, 1, $false, "[27]", 2) | Out-Null

# Reference list entry: "[22] https://interestingengineering.com/..." -> "[27]"
$r = $d.Paragraphs.Item(38).Range
$r.Find.Execute("[22]", $false, $false, $false, $false, $false, $true, 1, $false, "[27]", 2) | Out-Null

Write-Output "done"
